$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3547.75
$ws.Range("I40").Value = 3184.75
$ws.Range("J40").Value = 4273.75
$ws.Range("K40").Value = 3184.75
$ws.Range("L40").Value = 4273.75
$ws.Range("M40").Value = -3009.75
$ws.Range("N40").Value = -4623.75

$ws.Range("H137").Value = 3750.352
$ws.Range("I137").Value = 2349.3572
$ws.Range("K137").Value = 7048.071599999999
$ws.Range("M137").Value = -4498.071599999999

$ws.Range("H141").Value = 1026.8837
$ws.Range("I141").Value = 918.561
$ws.Range("K141").Value = 2755.683
$ws.Range("M141").Value = 2424.317

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4900.039
$ws.Range("I61").Value = 4138.5835
$ws.Range("J61").Value = 17083.334
$ws.Range("K61").Value = 4138.5835
$ws.Range("L61").Value = 17083.334
$ws.Range("M61").Value = -3926.5835
$ws.Range("N61").Value = -17507.334

$ws.Range("H74").Value = 246186.64
$ws.Range("I74").Value = 287391.78
$ws.Range("K74").Value = 287391.78
$ws.Range("M74").Value = -286517.78

$ws.Range("H77").Value = 246186.64
$ws.Range("I77").Value = 287391.78
$ws.Range("K77").Value = 1436958.9
$ws.Range("M77").Value = -1432590.9

$ws.Range("H97").Value = 1428585.2
$ws.Range("I97").Value = 1953443.2
$ws.Range("K97").Value = 1953443.2
$ws.Range("M97").Value = -1952947.2

$ws.Range("H132").Value = 3026.9688
$ws.Range("I132").Value = 2143.7625
$ws.Range("J132").Value = 7443
$ws.Range("K132").Value = 6431.287499999999
$ws.Range("L132").Value = 22329
$ws.Range("M132").Value = -3901.287499999999
$ws.Range("N132").Value = -27389

$ws.Range("H136").Value = 4900.039
$ws.Range("I136").Value = 4138.5835
$ws.Range("J136").Value = 17083.334
$ws.Range("K136").Value = 12415.7505
$ws.Range("L136").Value = 51250.00199999999
$ws.Range("M136").Value = -9865.750499999998
$ws.Range("N136").Value = -56350.00199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2750.1667
$ws.Range("I22").Value = 3500.3333
$ws.Range("K22").Value = 3500.3333
$ws.Range("M22").Value = -3327.3333

$ws.Range("H99").Value = 2406.6365
$ws.Range("I99").Value = 2285.4285
$ws.Range("K99").Value = 2285.4285
$ws.Range("M99").Value = -787.4285

$ws.Range("H105").Value = 40003316
$ws.Range("I105").Value = 111113800
$ws.Range("J105").Value = 3671.75
$ws.Range("K105").Value = 111113800
$ws.Range("L105").Value = 3671.75
$ws.Range("M105").Value = -111112053
$ws.Range("N105").Value = -7165.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4865.7803
$ws.Range("I31").Value = 2416.5
$ws.Range("J31").Value = 5655.871
$ws.Range("K31").Value = 2416.5
$ws.Range("L31").Value = 5655.871
$ws.Range("M31").Value = -2121.5
$ws.Range("N31").Value = -6245.871

$ws.Range("H34").Value = 4865.7803
$ws.Range("I34").Value = 2416.5
$ws.Range("J34").Value = 5655.871
$ws.Range("K34").Value = 2416.5
$ws.Range("L34").Value = 5655.871
$ws.Range("M34").Value = -2214.5
$ws.Range("N34").Value = -6059.871

$ws.Range("H58").Value = 3011.9092
$ws.Range("I58").Value = 2195.5366
$ws.Range("J58").Value = 5402.7144
$ws.Range("K58").Value = 2195.5366
$ws.Range("L58").Value = 5402.7144
$ws.Range("M58").Value = -1992.5366
$ws.Range("N58").Value = -5808.7144

$ws.Range("H107").Value = 412.83334
$ws.Range("I107").Value = 364.6
$ws.Range("J107").Value = 447.2857
$ws.Range("K107").Value = 364.6
$ws.Range("L107").Value = 447.2857
$ws.Range("M107").Value = 1555.4
$ws.Range("N107").Value = -4287.2857

$ws.Range("H136").Value = 3011.9092
$ws.Range("I136").Value = 2195.5366
$ws.Range("J136").Value = 5402.7144
$ws.Range("K136").Value = 6586.6098
$ws.Range("L136").Value = 16208.1432
$ws.Range("M136").Value = -4036.6098
$ws.Range("N136").Value = -21308.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5741.2256
$ws.Range("I131").Value = 991.125
$ws.Range("J131").Value = 7393.4346
$ws.Range("K131").Value = 2973.375
$ws.Range("L131").Value = 22180.3038
$ws.Range("M131").Value = 2066.625
$ws.Range("N131").Value = -32260.3038

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 35995
$ws.Range("I15").Value = 37000
$ws.Range("J15").Value = 34990
$ws.Range("K15").Value = 37000
$ws.Range("L15").Value = 34990
$ws.Range("M15").Value = -36712
$ws.Range("N15").Value = -35566

$ws.Range("H70").Value = 4997
$ws.Range("I70").Value = 4998
$ws.Range("J70").Value = 4994
$ws.Range("K70").Value = 4998
$ws.Range("L70").Value = 4994
$ws.Range("M70").Value = -4728
$ws.Range("N70").Value = -5534

$ws.Range("H73").Value = 4997
$ws.Range("I73").Value = 4998
$ws.Range("J73").Value = 4994
$ws.Range("K73").Value = 4998
$ws.Range("L73").Value = 4994
$ws.Range("M73").Value = -4062
$ws.Range("N73").Value = -6866

$ws.Range("H81").Value = 35995
$ws.Range("I81").Value = 37000
$ws.Range("J81").Value = 34990
$ws.Range("K81").Value = 37000
$ws.Range("L81").Value = 34990
$ws.Range("M81").Value = -36002
$ws.Range("N81").Value = -36986

$ws.Range("H84").Value = 35995
$ws.Range("I84").Value = 37000
$ws.Range("J84").Value = 34990
$ws.Range("K84").Value = 111000
$ws.Range("L84").Value = 104970
$ws.Range("M84").Value = -106008
$ws.Range("N84").Value = -114954

$ws.Range("H126").Value = 3374.4707
$ws.Range("I126").Value = 2346.8
$ws.Range("J126").Value = 4842.5713
$ws.Range("K126").Value = 7040.400000000001
$ws.Range("L126").Value = 14527.7139
$ws.Range("M126").Value = -4570.400000000001
$ws.Range("N126").Value = -19467.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2133.3333
$ws.Range("I82").Value = 1975
$ws.Range("J82").Value = 2450
$ws.Range("K82").Value = 1975
$ws.Range("L82").Value = 2450
$ws.Range("M82").Value = -1614
$ws.Range("N82").Value = -3172

$ws.Range("H85").Value = 2133.3333
$ws.Range("I85").Value = 1975
$ws.Range("J85").Value = 2450
$ws.Range("K85").Value = 1975
$ws.Range("L85").Value = 2450
$ws.Range("M85").Value = -727
$ws.Range("N85").Value = -4946

$ws.Range("H93").Value = 4759.769
$ws.Range("I93").Value = 4721.778
$ws.Range("J93").Value = 4845.25
$ws.Range("K93").Value = 4721.778
$ws.Range("L93").Value = 4845.25
$ws.Range("M93").Value = -3473.778
$ws.Range("N93").Value = -7341.25

$ws.Range("H136").Value = 2817.628
$ws.Range("I136").Value = 2462.3142
$ws.Range("J136").Value = 4372.125
$ws.Range("K136").Value = 7386.942599999999
$ws.Range("L136").Value = 13116.375
$ws.Range("M136").Value = -4836.942599999999
$ws.Range("N136").Value = -18216.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9509.514999999999
$ws.Range("I81").Value = 2855.6
$ws.Range("J81").Value = 14499.95
$ws.Range("K81").Value = 5711.2
$ws.Range("L81").Value = 28999.9
$ws.Range("M81").Value = -4650.2
$ws.Range("N81").Value = -31121.9

$ws.Range("H84").Value = 9509.514999999999
$ws.Range("I84").Value = 2855.6
$ws.Range("J84").Value = 14499.95
$ws.Range("K84").Value = 28556
$ws.Range("L84").Value = 144999.5
$ws.Range("M84").Value = -23252
$ws.Range("N84").Value = -155607.5

$ws.Range("H133").Value = 113887.5
$ws.Range("J133").Value = 113887.5
$ws.Range("L133").Value = 113887.5
$ws.Range("N133").Value = -124007.5

$ws.Range("H136").Value = 1514.1094
$ws.Range("I136").Value = 572.6226
$ws.Range("J136").Value = 6050.364
$ws.Range("K136").Value = 1717.8678
$ws.Range("L136").Value = 18151.092
$ws.Range("M136").Value = 832.1322
$ws.Range("N136").Value = -23251.092
